$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates are Excel serial numbers for 2021-05-03 .. 2021-05-05)
$newRows = @(
    @{ Row = 245; Date = 44319; B = 0; C = 12; D = 141.3927182750088 },
    @{ Row = 246; Date = 44320; B = 1; C = 9;  D = 106.0445387062566 },
    @{ Row = 247; Date = 44321; B = 1; C = 10; D = 117.827265229174 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.Date
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D

    # Copy the date style (number format + alignment + border) from the row above
    $ws.Cells.Item($rowNum - 1, 1).Copy()
    $ws.Cells.Item($rowNum, 1).PasteSpecial(-4122) | Out-Null
}
